# Apply the 2022-05-12 Fonds de solidarite volet 1 data update.
# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# for the affected rows of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 181347
$ws.Range("E8").Value = 650004878
$ws.Range("C10").Value = 278186
$ws.Range("D10").Value = 29588
$ws.Range("E10").Value = 1752125292
$ws.Range("C110").Value = 16867
$ws.Range("E110").Value = 25929535
$ws.Range("C115").Value = 17541
$ws.Range("E115").Value = 38583341
$ws.Range("C116").Value = 5253
$ws.Range("E116").Value = 15510059
$ws.Range("C117").Value = 19697
$ws.Range("E117").Value = 56402512
$ws.Range("C121").Value = 5961
$ws.Range("E121").Value = 11511893
$ws.Range("C131").Value = 7756
$ws.Range("E131").Value = 16748126
$ws.Range("C134").Value = 5665
$ws.Range("E134").Value = 17027467
$ws.Range("C138").Value = 2837
$ws.Range("E138").Value = 6576697
$ws.Range("C157").Value = 21202
$ws.Range("E157").Value = 77993905
$ws.Range("C161").Value = 44658
$ws.Range("E161").Value = 146110447
$ws.Range("C168").Value = 284906
$ws.Range("E168").Value = 1207681746
$ws.Range("C170").Value = 367228
$ws.Range("E170").Value = 2843315793
$ws.Range("C171").Value = 115091
$ws.Range("E171").Value = 444562942
$ws.Range("C174").Value = 357141
$ws.Range("D174").Value = 69786
$ws.Range("E174").Value = 1015883833
$ws.Range("C175").Value = 125491
$ws.Range("E175").Value = 810671097
$ws.Range("C179").Value = 235636
$ws.Range("E179").Value = 811898627
$ws.Range("C186").Value = 21930
$ws.Range("E186").Value = 40006767
$ws.Range("C188").Value = 19695
$ws.Range("E188").Value = 65984473
$ws.Range("C192").Value = 7457
$ws.Range("E192").Value = 17063682
$ws.Range("C193").Value = 5343
$ws.Range("E193").Value = 27693629
$ws.Range("C196").Value = 7396
$ws.Range("E196").Value = 20644396
$ws.Range("C203").Value = 13094
$ws.Range("E203").Value = 32950273
$ws.Range("C205").Value = 11114
$ws.Range("E205").Value = 43983889
$ws.Range("C213").Value = 3629
$ws.Range("E213").Value = 11065015
$ws.Range("C257").Value = 182547
$ws.Range("E257").Value = 1063782199
$ws.Range("C258").Value = 15139
$ws.Range("E258").Value = 40672103
$ws.Range("C266").Value = 71662
$ws.Range("E266").Value = 219427243
$ws.Range("C293").Value = 61659
$ws.Range("E293").Value = 194850617
$ws.Range("C295").Value = 91331
$ws.Range("E295").Value = 552907076
$ws.Range("C307").Value = 39618
$ws.Range("E307").Value = 95374311
$ws.Range("C313").Value = 220630
$ws.Range("E313").Value = 1370608294
$ws.Range("C322").Value = 81159
$ws.Range("E322").Value = 254507703
